# Apply cryptos list update (prices + 1h volumes) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.764.53"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "2.307.24"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'322.29"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "'105.23"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").Value = "'40.46"
$ws.Range("E10").Value = "  +4.41%  "
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "'8.61"
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'0.977"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "'15.37"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "2.655.53"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "2.298.90"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "42.816.74"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "'13.46"
$ws.Range("E21").Value = "  +35.04%  "
$ws.Range("D22").Value = "'73.82"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "'3.60"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "'272.19"
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").Value = "'2.25"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").Value = "'10.97"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "'22.74"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'38.29"
$ws.Range("E30").Value = "  +11.30%  "
$ws.Range("D31").Value = "'165.70"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "'6.18"
$ws.Range("E32").Value = "  +6.53%  "
$ws.Range("D33").Value = "'0.0889"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").Value = "'0.132"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.116"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.54"
$ws.Range("E36").Value = "  -12.22%  "
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "'0.0357"
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("E41").Value = "  +7.72%  "
$ws.Range("D42").Value = "'101.09"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "'70.76"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.226"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").Value = "'12.50"
$ws.Range("E46").Value = "  +6.12%  "
$ws.Range("D47").Value = "'82.84"
$ws.Range("E47").Value = "  +10.28%  "
$ws.Range("D48").Value = "'114.11"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.92"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.31"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "1.595.15"
$ws.Range("E51").Value = "  +4.72%  "
